# Super Bowl Week update
# Updates odds for rows 274-279 (week 19 games) and appends rows 280-286
# (weeks 20-22: divisional round, conference championships, Super Bowl).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows (only the odds columns D:G changed) ---

# Row 274: cincinnati-bengals vs las-vegas-raiders
$ws.Range("D274").Value = -280
$ws.Range("E274").Value = 225
$ws.Range("F274").Value = 6.5
$ws.Range("G274").Value = 48.5

# Row 275: buffalo-bills vs new-england-patriots
$ws.Range("D275").Value = -210
$ws.Range("E275").Value = 175
$ws.Range("F275").Value = 4
$ws.Range("G275").Value = 43

# Row 276: tampa-bay-buccaneers vs philadelphia-eagles
$ws.Range("D276").Value = -335
$ws.Range("E276").Value = 260
$ws.Range("F276").Value = 7.5
$ws.Range("G276").Value = 47.5

# Row 277: dallas-cowboys vs san-francisco-49ers
$ws.Range("D277").Value = -170
$ws.Range("E277").Value = 150
$ws.Range("F277").Value = 3.5
$ws.Range("G277").Value = 51

# Row 278: kansas-city-chiefs vs pittsburgh-steelers
$ws.Range("D278").Value = -600
$ws.Range("E278").Value = 435
$ws.Range("F278").Value = 11
$ws.Range("G278").Value = 46.5

# Row 279: los-angeles-rams vs arizona-cardinals
$ws.Range("D279").Value = -170
$ws.Range("E279").Value = 150

# --- Append new rows for weeks 20-22 (divisional, conference champ, Super Bowl) ---

$newRows = @(
    @(278, "tennessee-titans",      "cincinnati-bengals",   -200,  170, 4,   48,   20, 2021),
    @(279, "green-bay-packers",     "san-francisco-49ers",  -240,  195, 5.5, 47,   20, 2021),
    @(280, "tampa-bay-buccaneers",  "los-angeles-rams",     -140,  120, 2.5, 48,   20, 2021),
    @(281, "kansas-city-chiefs",    "buffalo-bills",        -130,  110, 2.5, 54,   20, 2021),
    @(282, "kansas-city-chiefs",    "cincinnati-bengals",   -350,  270, 7,   54.5, 21, 2021),
    @(283, "los-angeles-rams",      "san-francisco-49ers",  -180,  155, 3.5, 45.5, 21, 2021),
    @(284, "cincinnati-bengals",    "los-angeles-rams",      170, -200, 4,   48.5, 22, 2021)
)

$startRow = 280
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]

    # Column A carries the bold/bordered/centered style used throughout the
    # table (same as every other row's A cell) - copy it from the row above.
    $ws.Range("A" + ($r - 1)).Copy()
    $ws.Range("A" + $r).PasteSpecial(-4122)
}

$excel.CutCopyMode = 0
